$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AC, AD, AE for Wins / Losses / Ties
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match the bold/centered/bordered header style used by the existing
# header cells (e.g. AB1) by copying its format onto the new header cells.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# Season record values for every data row (2-43): 75 wins, 87 losses, 0 ties
$lastRow = 43
$ws.Range("AC2:AC$lastRow").Value = 75
$ws.Range("AD2:AD$lastRow").Value = 87
$ws.Range("AE2:AE$lastRow").Value = 0
